$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 37081.25
$ws.Range("J81").Value = 37081.25
$ws.Range("L81").Value = 37081.25
$ws.Range("N81").Value = -39077.25
$ws.Range("H84").Value = 37081.25
$ws.Range("J84").Value = 37081.25
$ws.Range("L84").Value = 111243.75
$ws.Range("N84").Value = -121227.75
$ws.Range("H87").Value = 29246
$ws.Range("J87").Value = 29246
$ws.Range("L87").Value = 29246
$ws.Range("N87").Value = -31742
$ws.Range("H90").Value = 29246
$ws.Range("J90").Value = 29246
$ws.Range("L90").Value = 87738
$ws.Range("N90").Value = -100218
$ws.Range("H94").Value = 6500
$ws.Range("I94").Value = 6500
$ws.Range("K94").Value = 6500
$ws.Range("M94").Value = -6049
$ws.Range("H115").Value = 1037.8572
$ws.Range("I115").Value = 475.45456
$ws.Range("J115").Value = 3100
$ws.Range("K115").Value = 1426.36368
$ws.Range("L115").Value = 9300
$ws.Range("M115").Value = 140.6363200000001
$ws.Range("N115").Value = -12434
$ws.Range("H132").Value = 16674976
$ws.Range("I132").Value = 20008020
$ws.Range("K132").Value = 60024060
$ws.Range("M132").Value = -60021530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13255.349
$ws.Range("I32").Value = 11913.926
$ws.Range("J32").Value = 17340.592
$ws.Range("K32").Value = 11913.926
$ws.Range("L32").Value = 17340.592
$ws.Range("M32").Value = -11626.926
$ws.Range("N32").Value = -17914.592
$ws.Range("H132").Value = 2751.85
$ws.Range("I132").Value = 2275.8462
$ws.Range("J132").Value = 3635.8572
$ws.Range("K132").Value = 6827.5386
$ws.Range("L132").Value = 10907.5716
$ws.Range("M132").Value = -4297.5386
$ws.Range("N132").Value = -15967.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 165.33333
$ws.Range("I22").Value = 165.33333
$ws.Range("K22").Value = 165.33333
$ws.Range("M22").Value = 7.666670000000011
$ws.Range("H69").Value = 30647.5
$ws.Range("J69").Value = 30647.5
$ws.Range("L69").Value = 30647.5
$ws.Range("N69").Value = -32269.5
$ws.Range("H72").Value = 30647.5
$ws.Range("J72").Value = 30647.5
$ws.Range("L72").Value = 91942.5
$ws.Range("N72").Value = -100054.5
$ws.Range("H99").Value = 3932.1538
$ws.Range("I99").Value = 3464.2727
$ws.Range("J99").Value = 6505.5
$ws.Range("K99").Value = 3464.2727
$ws.Range("L99").Value = 6505.5
$ws.Range("M99").Value = -1966.2727
$ws.Range("N99").Value = -9501.5
$ws.Range("H132").Value = 21286.215
$ws.Range("J132").Value = 21286.215
$ws.Range("L132").Value = 21286.215
$ws.Range("N132").Value = -31406.215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3993.1865
$ws.Range("I31").Value = 2978.6667
$ws.Range("J31").Value = 4849.1875
$ws.Range("K31").Value = 2978.6667
$ws.Range("L31").Value = 4849.1875
$ws.Range("M31").Value = -2683.6667
$ws.Range("N31").Value = -5439.1875
$ws.Range("H34").Value = 3993.1865
$ws.Range("I34").Value = 2978.6667
$ws.Range("J34").Value = 4849.1875
$ws.Range("K34").Value = 2978.6667
$ws.Range("L34").Value = 4849.1875
$ws.Range("M34").Value = -2776.6667
$ws.Range("N34").Value = -5253.1875
$ws.Range("H68").Value = 29936.8
$ws.Range("J68").Value = 29936.8
$ws.Range("L68").Value = 29936.8
$ws.Range("N68").Value = -31434.8
$ws.Range("H70").Value = 38114.285
$ws.Range("I70").Value = 37000
$ws.Range("J70").Value = 38300
$ws.Range("K70").Value = 37000
$ws.Range("L70").Value = 38300
$ws.Range("M70").Value = -36685
$ws.Range("N70").Value = -38930
$ws.Range("H71").Value = 29936.8
$ws.Range("J71").Value = 29936.8
$ws.Range("L71").Value = 89810.39999999999
$ws.Range("N71").Value = -97298.39999999999
$ws.Range("H73").Value = 38114.285
$ws.Range("I73").Value = 37000
$ws.Range("J73").Value = 38300
$ws.Range("K73").Value = 37000
$ws.Range("L73").Value = 38300
$ws.Range("M73").Value = -35908
$ws.Range("N73").Value = -40484
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H99").Value = 4967.75
$ws.Range("I99").Value = 2675
$ws.Range("J99").Value = 6114.125
$ws.Range("K99").Value = 2675
$ws.Range("L99").Value = 6114.125
$ws.Range("M99").Value = -1177
$ws.Range("N99").Value = -9110.125
$ws.Range("H122").Value = 2997.9565
$ws.Range("I122").Value = 2619.0588
$ws.Range("J122").Value = 4071.5
$ws.Range("K122").Value = 7857.176399999999
$ws.Range("L122").Value = 12214.5
$ws.Range("M122").Value = -5407.176399999999
$ws.Range("N122").Value = -17114.5
$ws.Range("H126").Value = 4967.75
$ws.Range("I126").Value = 2675
$ws.Range("J126").Value = 6114.125
$ws.Range("K126").Value = 8025
$ws.Range("L126").Value = 18342.375
$ws.Range("M126").Value = -5555
$ws.Range("N126").Value = -23282.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1668.6052
$ws.Range("I131").Value = 3345.875
$ws.Range("J131").Value = 1221.3334
$ws.Range("K131").Value = 10037.625
$ws.Range("L131").Value = 3664.0002
$ws.Range("M131").Value = -4997.625
$ws.Range("N131").Value = -13744.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 15000
$ws.Range("J86").Value = 15000
$ws.Range("L86").Value = 15000
$ws.Range("N86").Value = -17372
$ws.Range("H89").Value = 15000
$ws.Range("J89").Value = 15000
$ws.Range("L89").Value = 45000
$ws.Range("N89").Value = -56856
$ws.Range("H102").Value = 34632.03
$ws.Range("I102").Value = 1868.8462
$ws.Range("K102").Value = 1868.8462
$ws.Range("M102").Value = -246.8462
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178
$ws.Range("H132").Value = 4077.8538
$ws.Range("I132").Value = 4400.227
$ws.Range("J132").Value = 3704.5789
$ws.Range("K132").Value = 13200.681
$ws.Range("L132").Value = 11113.7367
$ws.Range("M132").Value = -10670.681
$ws.Range("N132").Value = -16173.7367

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2228.5
$ws.Range("I93").Value = 1655.4445
$ws.Range("J93").Value = 3260
$ws.Range("K93").Value = 1655.4445
$ws.Range("L93").Value = 3260
$ws.Range("M93").Value = -407.4445000000001
$ws.Range("N93").Value = -5756

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 9807.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 9807.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 9807.75
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = -10827.75
$ws.Range("H55").Value = 7397.4
$ws.Range("J55").Value = 8496.75
$ws.Range("L55").Value = 8496.75
$ws.Range("N55").Value = -9050.75
$ws.Range("H59").Value = 15000
$ws.Range("J59").Value = 15000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16476
$ws.Range("H61").Value = 21350
$ws.Range("I61").Value = 4050
$ws.Range("K61").Value = 4050
$ws.Range("M61").Value = -3758
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360
$ws.Range("H122").Value = 2034.35
$ws.Range("I122").Value = 1628.6471
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 4885.9413
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -2435.9413
$ws.Range("N122").Value = -17900.0005
